$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.879.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.536.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.585'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.95%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.547'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.11'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0816'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.59'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.921.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.570.09'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.864'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.890.62'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0984'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.10'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +10.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.10'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '158.41'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.49%  '
$ws.Range("E33").Value = '  -1.66%  '
$ws.Range("E34").Value = '  -2.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0795'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.41%  '
$ws.Range("E36").Value = '  -4.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.117'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.39'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("E39").Value = '  +0.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.32'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.64%  '
$ws.Range("E42").Value = '  -1.14%  '
$ws.Range("E43").Value = '  +0.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0304'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.051.12'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '86.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.98'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.781.14'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.192'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '103.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.28%  '
